$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (A5): remove the word "proper " before "usage of personal protective equipment"
$ws.Range("A5").Value = "6. Strategies to reduce transmission and acquisition in healthcare settings. This section encompasses approaches such as education, emphasizing risk assessment, and  usage of personal protective equipment (PPE, such as gloves, aprons, face masks) for cystic fibrosis (CF) patients and / or healthcare professionals both within hospital premises and in the community, and the involvement of family and friends. It considers hand hygiene and  environmental cleaning and disinfection within hospitals.  It includes  water safety, ventilation quality, disinfection protocols for rooms and equipment, and specific considerations for pulmonary function testing areas involving negative pressure rooms, HEPA filters, timing between patients, and UV germicidal irradiation.  Audit of infection control is also considered."

# Row 6 (A6): replace "7. Strategies for healthcare settings..." topic text with the new
# "7. Strategies for non-healthcare settings..." topic text
$ws.Range("A6").Value = "7. Strategies for non-healthcare settings. This section covers various domains including home, indoor settings, gyms, outdoor environments, educational settings from pre-school to higher education, and travel considerations. It emphasizes strategies like cleaning rooms and equipment, ensuring water safety, addressing potential risks in outdoor activities, and implementing preventive measures in educational settings and during travel, especially air travel. Additionally, it highlights the importance of viral transmission prevention and immunization against diseases such as flu, COVID-19, RSV, and pneumococcus, aligning with considerations in the CF Trust Antibiotic Treatment guideline."
